$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.900.07"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "3.846.85"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "704.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").Value = "3.844.66"
$ws.Range("E7").Value = "  +1.03%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("E11").Value = "  -1.73%  "
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("E13").Value = "  -2.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").Value = "4.496.61"
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("D16").Value = "3.937.02"
$ws.Range("E16").Value = "  +3.49%  "
$ws.Range("D17").Value = "71.005.76"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "492.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.09%  "
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  -1.18%  "
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.91%  "
$ws.Range("D35").Value = "3.803.20"
$ws.Range("E35").Value = "  +1.26%  "
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E41").Value = "  +6.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.86%  "
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "163.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("E46").Value = "  -6.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "48.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.66%  "
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "412.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.13%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.84%  "
